$wb = $excel.ActiveWorkbook

# --- Update product names in the "sort" sheet's shared strings ---
$sortWs = $wb.Worksheets.Item("sort")
$sortWs.Range("C4").Value = "Sauce Labs Onesie"
$sortWs.Range("C5").Value = "Sauce Labs Fleece Jacket"

# --- Restyle the data range to Times New Roman ---
$dataRng = $sortWs.Range("A1:D5")
$dataRng.Font.Color = 0
$dataRng.Font.Name = "Times New Roman"
$dataRng.Font.Family = 1

# --- Switch the active tab from "login" to "sort" ---
$sortWs.Activate()
$sortWs.Range("C12").Select()
